# Re-apply the "Medium Style 2 - Accent 1" built-in table style to the
# financial-documents table on slide 5 (Google Shape;122;p17), replacing the
# presentation's custom Table_0 style ({DC0CE955-88B5-45CD-B694-8E140BC90E2F})
# with the built-in gallery style {43DACE84-D1A9-44B1-9E08-886F8A9F36A8}, as
# would happen by picking a different style in the Table Design > Table
# Styles gallery.

$p = $ppt.ActivePresentation

$targetStyleId = "{43DACE84-D1A9-44B1-9E08-886F8A9F36A8}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
